$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.938999999999991
$ws.Range("A10").Value = -22.05099999999999
$ws.Range("A12").Value = -21.48679999999998
$ws.Range("B12").Value = 5.494099999999996
$ws.Range("D12").Value = -5.888099999999998
$ws.Range("D13").Value = -8.658399999999984
$ws.Range("B17").Value = 5.326099999999997
$ws.Range("A18").Value = -22.1941
$ws.Range("D21").Value = -8.135199999999996
$ws.Range("B26").Value = 4.093900000000002
$ws.Range("B27").Value = 6.720200000000005
$ws.Range("B28").Value = 6.442299999999999
$ws.Range("D36").Value = -7.210800000000003
$ws.Range("A37").Value = -21.80429999999999
$ws.Range("B37").Value = 7.021699999999999
$ws.Range("D38").Value = -8.0733
$ws.Range("D41").Value = -8.217400000000003
$ws.Range("D52").Value = -7.834900000000006
$ws.Range("A55").Value = -22.0551
$ws.Range("B65").Value = 6.2351
$ws.Range("D67").Value = -7.474599999999997
$ws.Range("A68").Value = -21.47959999999999
$ws.Range("B73").Value = 8.811399999999997
$ws.Range("A77").Value = -19.91539999999999
$ws.Range("A78").Value = -19.89669999999997
$ws.Range("B84").Value = 5.772400000000001
$ws.Range("B85").Value = 6.004599999999998
$ws.Range("D89").Value = -8.193599999999998
$ws.Range("B93").Value = 5.529999999999997
$ws.Range("B95").Value = 5.143400000000005
$ws.Range("D95").Value = -7.569600000000003
$ws.Range("B98").Value = 4.866000000000009
$ws.Range("B99").Value = 5.966100000000001
$ws.Range("B101").Value = 5.992899999999999
$ws.Range("D105").Value = -8.104700000000005
